# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (plain Office palette)
#   ppt/theme/theme2.xml -> "Integral"     (the theme actually applied
#                                           to the slide master/deck)
#
# The authored commit swaps the two themes' contents, so the theme
# that is actually applied to the presentation switches from
# "Integral" to the plain "Office Theme" palette (and vice versa for
# the otherwise-unused second theme part).
#
# The PowerPoint object model doesn't expose a way to rename a theme
# or swap which XML part backs it, but it does expose the twelve
# theme colours themselves - dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink, in that fixed order - through
# Master.Theme.ThemeColorScheme(n).RGB. Driving every slot to the
# "Office Theme" RGB values reproduces the swap's visible effect on
# the theme that is actually applied to the presentation.

function ToRGB([byte]$r, [byte]$g, [byte]$b) {
    return [int]$r + ([int]$g * 256) + ([int]$b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme colour scheme values (previously ppt/theme/theme1.xml),
# applied in clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    (ToRGB 0x00 0x00 0x00),  # 1  dk1
    (ToRGB 0xFF 0xFF 0xFF),  # 2  lt1
    (ToRGB 0x44 0x54 0x6A),  # 3  dk2
    (ToRGB 0xE7 0xE6 0xE6),  # 4  lt2
    (ToRGB 0x5B 0x9B 0xD5),  # 5  accent1
    (ToRGB 0xED 0x7D 0x31),  # 6  accent2
    (ToRGB 0xA5 0xA5 0xA5),  # 7  accent3
    (ToRGB 0xFF 0xC0 0x00),  # 8  accent4
    (ToRGB 0x44 0x72 0xC4),  # 9  accent5
    (ToRGB 0x70 0xAD 0x47),  # 10 accent6
    (ToRGB 0x05 0x63 0xC1),  # 11 hlink
    (ToRGB 0x95 0x4F 0x72)   # 12 folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
